# Use rxridge With STATA 18
#
# - Rename the single "Regression Panel Data In R Olah" sheet to "Sheet1"
# - Add two more blank sheets, "Sheet2" and "Sheet3", right after it
# - Re-activate Sheet1 and select the whole grid on it (A1:XFD1048576)
# - Resize/reposition the workbook window to match the new saved view

$wb = $excel.ActiveWorkbook

# Rename the original (only) worksheet to "Sheet1"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# Insert two fresh worksheets right after Sheet1, named Sheet2 and Sheet3
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# Leave Sheet1 as the active/selected sheet, with its entire grid selected
$ws1.Activate() | Out-Null
$ws1.Cells.Select() | Out-Null

# Match the workbook window size/position recorded on save
$win = $excel.ActiveWindow
[void]($win.Left = 480)
[void]($win.Top = 360)
[void]($win.Width = 19815)
[void]($win.Height = 7650)
